$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 259 (shifts existing rows 259..332 down to 260..333)
$ws.Rows.Item(259).Insert()

# Populate the new row 259 with the data for the new weekly record
$ws.Cells.Item(259, 1).Value = 10
$ws.Cells.Item(259, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(259, 3).Value = "La Araucanía"
$ws.Cells.Item(259, 4).Value = 45173
$ws.Cells.Item(259, 5).Value = 9
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100104
$ws.Cells.Item(259, 8).Value = "Frutos de pepita"
$ws.Cells.Item(259, 9).Value = 100104003
$ws.Cells.Item(259, 10).Value = "Membrillo"
$ws.Cells.Item(259, 11).Value = "Champion"
$ws.Cells.Item(259, 12).Value = "Primera"
$ws.Cells.Item(259, 13).Value = 150
$ws.Cells.Item(259, 14).Value = 16000
$ws.Cells.Item(259, 15).Value = 16000
$ws.Cells.Item(259, 16).Value = 16000
$ws.Cells.Item(259, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(259, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(259, 19).Value = 889
$ws.Cells.Item(259, 20).Value = 18
